$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I11").Value = 336.66666
$ws.Range("K11").Value = 336.66666
$ws.Range("H11").Value = 336.66666
$ws.Range("M11").Value = -196.66666
$ws.Range("J32").Value = 6521.2
$ws.Range("N32").Value = -7173.2
$ws.Range("L32").Value = 6521.2
$ws.Range("H32").Value = 7828.2
$ws.Range("K33").Value = 476.13333
$ws.Range("M33").Value = -247.13333
$ws.Range("I33").Value = 476.13333
$ws.Range("H33").Value = 477.2
$ws.Range("N33").Value = -938.4
$ws.Range("J33").Value = 480.4
$ws.Range("L33").Value = 480.4
$ws.Range("K106").Value = 2290.923
$ws.Range("L106").Value = 6660.5
$ws.Range("H106").Value = 3670.7896
$ws.Range("N106").Value = -7922.5
$ws.Range("J106").Value = 6660.5
$ws.Range("I106").Value = 2290.923
$ws.Range("M106").Value = -1659.923
$ws.Range("M111").Value = -504.0001999999999
$ws.Range("J111").Value = 1200
$ws.Range("N111").Value = -9734
$ws.Range("I111").Value = 1190.3334
$ws.Range("L111").Value = 3600
$ws.Range("K111").Value = 3571.0002
$ws.Range("H111").Value = 1192.75
$ws.Range("M116").Value = -14556.2
$ws.Range("I116").Value = 17998.2
$ws.Range("H116").Value = 18187.125
$ws.Range("K116").Value = 17998.2
$ws.Range("J134").Value = 66580.75
$ws.Range("N134").Value = -76720.75
$ws.Range("H134").Value = 66580.75
$ws.Range("L134").Value = 66580.75
$ws.Range("H138").Value = 2517.104
$ws.Range("L138").Value = 10147.404
$ws.Range("I138").Value = 1161.3667
$ws.Range("N138").Value = -20427.404
$ws.Range("M138").Value = 1655.8999
$ws.Range("J138").Value = 3382.468
$ws.Range("K138").Value = 3484.1001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 5484.913
$ws.Range("I32").Value = 5484.913
$ws.Range("J32").Value = 10427.571
$ws.Range("N32").Value = -11001.571
$ws.Range("L32").Value = 10427.571
$ws.Range("M32").Value = -5197.913
$ws.Range("H32").Value = 6638.2
$ws.Range("N45").Value = -2309.3334
$ws.Range("J45").Value = 1555.3334
$ws.Range("K45").Value = 1616.9131
$ws.Range("I45").Value = 1616.9131
$ws.Range("M45").Value = -1239.9131
$ws.Range("H45").Value = 1609.8077
$ws.Range("L45").Value = 1555.3334
$ws.Range("L97").Value = 729
$ws.Range("H97").Value = 741.92
$ws.Range("N97").Value = -1721
$ws.Range("J97").Value = 729

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3059.2
$ws.Range("N99").Value = -6246
$ws.Range("J99").Value = 3250
$ws.Range("L99").Value = 3250
$ws.Range("L105").Value = 11674.75
$ws.Range("M105").Value = -24964.5
$ws.Range("I105").Value = 26711.5
$ws.Range("J105").Value = 11674.75
$ws.Range("K105").Value = 26711.5
$ws.Range("N105").Value = -15168.75
$ws.Range("H105").Value = 19193.125
$ws.Range("I134").Value = 1104.6842
$ws.Range("H134").Value = 1237.238
$ws.Range("K134").Value = 3314.0526
$ws.Range("M134").Value = -779.0526

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J62").Value = 16202.25
$ws.Range("I62").Value = 3743.5
$ws.Range("H62").Value = 9972.875
$ws.Range("L62").Value = 16202.25
$ws.Range("K62").Value = 3743.5
$ws.Range("M62").Value = -3119.5
$ws.Range("N62").Value = -17450.25
$ws.Range("L65").Value = 81011.25
$ws.Range("K65").Value = 18717.5
$ws.Range("I65").Value = 3743.5
$ws.Range("H65").Value = 9972.875
$ws.Range("J65").Value = 16202.25
$ws.Range("N65").Value = -87251.25
$ws.Range("M65").Value = -15597.5
$ws.Range("J103").Value = 0
$ws.Range("H103").Value = 49986.668
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -48814.668
$ws.Range("N103").ClearContents()
$ws.Range("I103").Value = 49986.668
$ws.Range("K103").Value = 49986.668
$ws.Range("L105").Value = 6488.75
$ws.Range("M105").Value = 765.55554
$ws.Range("I105").Value = 981.44446
$ws.Range("J105").Value = 6488.75
$ws.Range("K105").Value = 981.44446
$ws.Range("N105").Value = -9982.75
$ws.Range("H105").Value = 2676
$ws.Range("K107").Value = 1735.9166
$ws.Range("L107").Value = 1821.6666
$ws.Range("I107").Value = 1735.9166
$ws.Range("J107").Value = 1821.6666
$ws.Range("H107").Value = 1764.5
$ws.Range("N107").Value = -5661.6666
$ws.Range("M107").Value = 184.0834

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 227083.83
$ws.Range("L37").Value = 681251.49
$ws.Range("N37").Value = -681475.49
$ws.Range("J37").Value = 227083.83
$ws.Range("M86").Value = -4518.9998
$ws.Range("K86").Value = 5704.9998
$ws.Range("H86").Value = 2642.25
$ws.Range("J86").Value = 3382.8333
$ws.Range("L86").Value = 10148.4999
$ws.Range("I86").Value = 1901.6666
$ws.Range("N86").Value = -12520.4999
$ws.Range("H89").Value = 2642.25
$ws.Range("K89").Value = 17114.9994
$ws.Range("M89").Value = -11186.9994
$ws.Range("I89").Value = 1901.6666
$ws.Range("J89").Value = 3382.8333
$ws.Range("L89").Value = 30445.4997
$ws.Range("N89").Value = -42301.4997
$ws.Range("K107").Value = 7535.400000000001
$ws.Range("L107").Value = 8526006.600000001
$ws.Range("I107").Value = 2511.8
$ws.Range("J107").Value = 2842002.2
$ws.Range("H107").Value = 1954661.5
$ws.Range("N107").Value = -8529846.600000001
$ws.Range("M107").Value = -5615.400000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J103").Value = 23931.2
$ws.Range("H103").Value = 23931.2
$ws.Range("L103").Value = 23931.2
$ws.Range("N103").Value = -26275.2
$ws.Range("J111").Value = 43731
$ws.Range("N111").Value = -49865
$ws.Range("L111").Value = 43731
$ws.Range("H111").Value = 43731
$ws.Range("M113").Value = 693.2858000000001
$ws.Range("J113").Value = 3597
$ws.Range("I113").Value = 1476.7142
$ws.Range("K113").Value = 1476.7142
$ws.Range("N113").Value = -7937
$ws.Range("H113").Value = 2112.8
$ws.Range("L113").Value = 3597
$ws.Range("I132").Value = 5825.7407
$ws.Range("L132").Value = 22202.571
$ws.Range("H132").Value = 6150.0293
$ws.Range("J132").Value = 7400.857
$ws.Range("N132").Value = -27262.571
$ws.Range("M132").Value = -14947.2221
$ws.Range("K132").Value = 17477.2221
$ws.Range("J134").Value = 55220.2
$ws.Range("N134").Value = -170730.6
$ws.Range("H134").Value = 55220.2
$ws.Range("L134").Value = 165660.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4153
$ws.Range("J46").Value = 4315.4614
$ws.Range("L46").Value = 4315.4614
$ws.Range("N46").Value = -4691.4614
$ws.Range("L55").Value = 4986
$ws.Range("I55").Value = 3571998.5
$ws.Range("K55").Value = 3571998.5
$ws.Range("N55").Value = -5332
$ws.Range("J55").Value = 4986
$ws.Range("H55").Value = 2382994.5
$ws.Range("M55").Value = -3571825.5
$ws.Range("K61").Value = 2981.4666
$ws.Range("H61").Value = 4133.273
$ws.Range("M61").Value = -2779.4666
$ws.Range("L61").Value = 6601.4287
$ws.Range("J61").Value = 6601.4287
$ws.Range("I61").Value = 2981.4666
$ws.Range("N61").Value = -7005.4287
$ws.Range("M113").Value = -811.4666000000002
$ws.Range("J113").Value = 6601.4287
$ws.Range("I113").Value = 2981.4666
$ws.Range("K113").Value = 2981.4666
$ws.Range("N113").Value = -10941.4287
$ws.Range("H113").Value = 4133.273
$ws.Range("L113").Value = 6601.4287
$ws.Range("I132").Value = 15814.667
$ws.Range("H132").Value = 14561.692
$ws.Range("M132").Value = -44914.001
$ws.Range("K132").Value = 47444.001
$ws.Range("L136").Value = 33043.749
$ws.Range("J136").Value = 11014.583
$ws.Range("H136").Value = 8048.7
$ws.Range("N136").Value = -38143.749

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J122").Value = 8501.666999999999
$ws.Range("I122").Value = 10555.5
$ws.Range("H122").Value = 9323.200000000001
$ws.Range("L122").Value = 25505.001
$ws.Range("K122").Value = 31666.5
$ws.Range("M122").Value = -29216.5
$ws.Range("N122").Value = -30405.001
